# "Actualización desde MV -datos-"
# Append the new quarterly data row (period 01-07-2021) to the bottom of
# the table on Sheet1, following the same layout as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 76

# Column A holds a date-like label ("01-07-2021") that must be stored as
# literal text (a shared string), exactly like all the other period
# labels above it in column A - not auto-converted into a date serial
# number. Temporarily force a Text number format before assigning the
# value, then restore the default format so the cell ends up with no
# explicit style, matching its neighbours.
$colA = $ws.Cells.Item($row, 1)
$colA.NumberFormat = "@"
$colA.Value = "01-07-2021"
$colA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 233155
$ws.Cells.Item($row, 3).Value = 2837
$ws.Cells.Item($row, 4).Value = 244
$ws.Cells.Item($row, 5).Value = 2593
$ws.Cells.Item($row, 6).Value = 5012
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 5011
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 36094
$ws.Cells.Item($row, 12).Value = 9803
$ws.Cells.Item($row, 13).Value = 26291
$ws.Cells.Item($row, 14).Value = 73426
$ws.Cells.Item($row, 15).Value = 13927
$ws.Cells.Item($row, 16).Value = 59500
$ws.Cells.Item($row, 17).Value = 112279
$ws.Cells.Item($row, 18).Value = 1428
$ws.Cells.Item($row, 19).Value = 110851
$ws.Cells.Item($row, 20).Value = 3506
$ws.Cells.Item($row, 21).Value = 3506
